# Populate the 'Bugs' worksheet (sheet2) per the commit's bug-report table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bugs")

# Column widths (A, B, C) matching the authored layout.
$ws.Columns.Item(1).ColumnWidth = 67.16666666666667
$ws.Columns.Item(2).ColumnWidth = 126.83333333333333
$ws.Columns.Item(3).ColumnWidth = 68.33333333333333

# --- Header row (row 2): bold labels A2:E2, trailing formatted-but-empty
# cells F2:Z2 left over from the source sheet's header-row paint range.
$ws.Range("A2").Value = "Título descriptivo"
$ws.Range("B2").Value = "Pasos para reproducir"
$ws.Range("C2").Value = "Resultado esperado vs obtenido"
$ws.Range("D2").Value = "Severidad"
$ws.Range("E2").Value = "Evidencia visual (screenshot o video corto)"
$headerRange = $ws.Range("A2:E2")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 11
$headerRange.Font.Color = 0
$headerRange.Font.Name = "Arial"

$headerTailRange = $ws.Range("F2:Z2")
$headerTailRange.Font.Bold = $true
$headerTailRange.Font.ThemeColor = 1
$headerTailRange.Font.Name = "Arial"

# --- Row 3: first bug (plain labels, B3 carries the light fill band).
$ws.Range("A3").Value = "Subir foto de la galería, muestra banner donde no se ve texto."
$ws.Range("B3").Value = "1. Click en el icono `"+`". 2. Elegir foto. 3. Publicar."
$ws.Range("C3").Value = "Se sube la imágen pero hay un error en el banner que dice siguiente, pues no se ve el texto."
$ws.Range("D3").Value = "Baja"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.Color = 0
$ws.Range("A3").Font.Name = "Arial"

$ws.Range("B3").Interior.Color = 16777215
$ws.Range("B3").Font.Size = 11
$ws.Range("B3").Font.Color = 2039583
$ws.Range("B3").Font.Name = "Arial"

$row3CD = $ws.Range("C3:D3")
$row3CD.Font.Size = 11
$row3CD.Font.ThemeColor = 1
$row3CD.Font.Name = "Arial"

$row3Tail = $ws.Range("E3:Z3")
$row3Tail.Font.Size = 11
$row3Tail.Font.ThemeColor = 1
$row3Tail.Font.Name = "Arial"

# --- Row 4: second bug (plain labels, same theme-coloured font as C3:D3).
$ws.Range("A4").Value = "Cierre forzado (Crash) al intentar subir imagen con menos de 10MB de espacio disponible."
$ws.Range("B4").Value = "1. Abrir Instagram Lite. 2. Tocar el icono `"+`" para crear una nueva publicación. 3. Seleccionar cualquier imagen de la galería. 4. Tocar en `"Siguiente`" y luego en `"Compartir`"."
$ws.Range("C4").Value = "La aplicación debe mostrar un mensaje de error amigable: `"No hay suficiente espacio en el dispositivo para procesar la imagen`" y permitir al usuario regresar al feed. En cambio la pantalla se queda en blanco por 3 segundos y aparece el mensaje de sistema: `"Instagram Lite se detuvo`"."
$ws.Range("D4").Value = "Crítica (El usuario no puede usar la función principal y la app se rompe)."
$row4 = $ws.Range("A4:D4")
$row4.Font.Size = 11
$row4.Font.ThemeColor = 1
$row4.Font.Name = "Arial"

$row4Tail = $ws.Range("E4:Z4")
$row4Tail.Font.Size = 11
$row4Tail.Font.ThemeColor = 1
$row4Tail.Font.Name = "Arial"

# --- Rows 5-7: remaining bugs, white fill band with wrapped description columns.
# Row 5 keeps Arial throughout (title/expected-vs-obtained col wraps);
# rows 6-7 switch the whole band to "Google Sans" (title col no longer wraps).
$ws.Range("A5").Value = "Bucle de carga infinita en la pestaña de Reels"
$ws.Range("B5").Value = "1. Abrir la aplicación Instagram Lite. 2. Hacer clic en el icono central de Reels. 3. Deslizar hacia abajo rápidamente para pasar 5 o 6 videos."
$ws.Range("C5").Value = " - Esperado: El siguiente video debe precargarse y reproducirse fluidamente.`n - Obtenido: La pantalla se queda en negro con un círculo de carga infinito, incluso con una conexión Wi-Fi estable."
$ws.Range("D5").Value = "Media/Alta (Afecta el consumo de contenido principal)"
$row5 = $ws.Range("A5:D5")
$row5.Interior.Color = 16777215
$row5.Font.Size = 11
$row5.Font.Color = 657930
$row5.Font.Name = "Arial"
$ws.Range("C5").WrapText = $true
$ws.Range("B5").Font.Name = "Google Sans"
$ws.Range("B5").WrapText = $true
$ws.Range("D5").Font.Name = "Google Sans"

$ws.Range("A6").Value = "Error de `"Imagen no disponible`" en Carruseles"
$ws.Range("B6").Value = "1. Seleccionar una publicación que contenga un carrusel (múltiples fotos). 2. Deslizar hacia la derecha para ver la segunda o tercera imagen. 3. Regresar rápidamente a la primera imagen."
$ws.Range("C6").Value = " - Esperado: Todas las imágenes del carrusel deben ser visibles al deslizar.`n - Obtenido: Aparece un icono de `"triángulo de advertencia`" o un cuadro gris, indicando que la imagen no se pudo cargar."
$ws.Range("D6").Value = "Media (Impacta la experiencia visual pero no bloquea la app)."
$rowAll = $ws.Range("A6:D6")
$rowAll.Interior.Color = 16777215
$rowAll.Font.Size = 11
$rowAll.Font.Color = 657930
$rowAll.Font.Name = "Google Sans"
$ws.Range("B6:C6").WrapText = $true

$ws.Range("A7").Value = "Cierre inesperado (Crash) al subir Stories con stickers"
$ws.Range("B7").Value = "1. Tocar el icono de `"+`" para crear una Story. 2. Tomar una foto o subir una de la galería. 3. Intentar añadir un sticker interactivo (Encuesta, Enlace o Música)."
$ws.Range("C7").Value = " - Esperado: El menú de stickers debe abrirse y permitir la selección.`n - Obtenido: La aplicación se congela por 2 segundos y se cierra inesperadamente, regresando al menú de inicio del teléfono."
$ws.Range("D7").Value = "Crítica (Interrumpe una función principal de creación)."
$rowAll = $ws.Range("A7:D7")
$rowAll.Interior.Color = 16777215
$rowAll.Font.Size = 11
$rowAll.Font.Color = 657930
$rowAll.Font.Name = "Google Sans"
$ws.Range("B7:C7").WrapText = $true

# --- Rows 8-10: trailing formatted-but-empty cells left over from the
# source sheet's fill-down range.
$ws.Range("B8:C8").Interior.Color = 16777215
$ws.Range("B8:C8").Font.Size = 11
$ws.Range("B8:C8").Font.Color = 657930
$ws.Range("B8:C8").Font.Name = "Google Sans"
$ws.Range("B8:C8").WrapText = $true

$ws.Range("A9").Font.Size = 11
$ws.Range("A9").Font.Color = 0
$ws.Range("A9").Font.Name = "Arial"

$ws.Range("B9").Interior.Color = 16777215
$ws.Range("B9").Font.Size = 11
$ws.Range("B9").Font.Color = 657930
$ws.Range("B9").Font.Name = "Google Sans"
$ws.Range("B9").WrapText = $true

$ws.Range("B10").Font.Bold = $true
$ws.Range("B10").Font.ThemeColor = 1
$ws.Range("B10").Font.Name = "Arial"

Write-Host "Bugs sheet populated."
